# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.863.85'
$ws.Range("E2").Value = '  +1.31%  '

$ws.Range("D3").Value = '3.471.85'
$ws.Range("E3").Value = '  +1.17%  '

$ws.Range("E4").Value = '  +0.02%  '

$__fmt = $ws.Range("D5").NumberFormat
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '414.26'
$ws.Range("D5").NumberFormat = $__fmt
$ws.Range("E5").Value = '  +1.09%  '

$__fmt = $ws.Range("D6").NumberFormat
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.19'
$ws.Range("D6").NumberFormat = $__fmt
$ws.Range("E6").Value = '  +0.46%  '

$ws.Range("E7").Value = '  -1.07%  '

$__fmt = $ws.Range("D8").NumberFormat
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").NumberFormat = $__fmt
$ws.Range("E8").Value = '  +0.02%  '

$__fmt = $ws.Range("D9").NumberFormat
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.726'
$ws.Range("D9").NumberFormat = $__fmt
$ws.Range("E9").Value = '  -2.11%  '

$__fmt = $ws.Range("D10").NumberFormat
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.148'
$ws.Range("D10").NumberFormat = $__fmt
$ws.Range("E10").Value = '  +4.03%  '

$__fmt = $ws.Range("D11").NumberFormat
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.58'
$ws.Range("D11").NumberFormat = $__fmt
$ws.Range("E11").Value = '  -0.73%  '

$__fmt = $ws.Range("D12").NumberFormat
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.64'
$ws.Range("D12").NumberFormat = $__fmt
$ws.Range("E12").Value = '  +4.38%  '

$__fmt = $ws.Range("D13").NumberFormat
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000220'
$ws.Range("D13").NumberFormat = $__fmt
$ws.Range("E13").Value = '  -2.15%  '

$ws.Range("D14").Value = '4.020.52'
$ws.Range("E14").Value = '  +1.29%  '

$__fmt = $ws.Range("D15").NumberFormat
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.141'
$ws.Range("D15").NumberFormat = $__fmt
$ws.Range("E15").Value = '  -0.21%  '

$ws.Range("E16").Value = '  -4.35%  '

$ws.Range("D17").Value = '3.453.25'
$ws.Range("E17").Value = '  -0.53%  '

$__fmt = $ws.Range("D18").NumberFormat
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.70'
$ws.Range("D18").NumberFormat = $__fmt
$ws.Range("E18").Value = '  +0.62%  '

$ws.Range("E19").Value = '  -1.61%  '

$ws.Range("D20").Value = '62.739.71'
$ws.Range("E20").Value = '  +1.10%  '

$__fmt = $ws.Range("D21").NumberFormat
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '466.09'
$ws.Range("D21").NumberFormat = $__fmt
$ws.Range("E21").Value = '  +3.18%  '

$__fmt = $ws.Range("D22").NumberFormat
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '91.10'
$ws.Range("D22").NumberFormat = $__fmt
$ws.Range("E22").Value = '  -0.80%  '

$ws.Range("E23").Value = '  +1.40%  '

$__fmt = $ws.Range("D24").NumberFormat
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.27'
$ws.Range("D24").NumberFormat = $__fmt
$ws.Range("E24").Value = '  +1.83%  '

$__fmt = $ws.Range("D25").NumberFormat
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.52'
$ws.Range("D25").NumberFormat = $__fmt
$ws.Range("E25").Value = '  +17.93%  '

$__fmt = $ws.Range("D26").NumberFormat
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.30'
$ws.Range("D26").NumberFormat = $__fmt
$ws.Range("E26").Value = '  +1.42%  '

$__fmt = $ws.Range("D27").NumberFormat
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '33.39'
$ws.Range("D27").NumberFormat = $__fmt
$ws.Range("E27").Value = '  +1.15%  '

$ws.Range("E28").Value = '  +0.46%  '

$ws.Range("E29").Value = '  -1.85%  '

$__fmt = $ws.Range("D30").NumberFormat
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.01'
$ws.Range("D30").NumberFormat = $__fmt
$ws.Range("E30").Value = '  -0.07%  '

$ws.Range("E31").Value = '  -2.54%  '

$ws.Range("E32").Value = '  -2.42%  '

$ws.Range("E33").Value = '  -2.12%  '

$__fmt = $ws.Range("D34").NumberFormat
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '40.74'
$ws.Range("D34").NumberFormat = $__fmt
$ws.Range("E34").Value = '  -5.72%  '

$ws.Range("E35").Value = '  +0.08%  '

$__fmt = $ws.Range("D36").NumberFormat
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '58.46'
$ws.Range("D36").NumberFormat = $__fmt
$ws.Range("E36").Value = '  +7.49%  '

$ws.Range("E37").Value = '  -2.73%  '

$ws.Range("E38").Value = '  +4.89%  '

$__fmt = $ws.Range("D39").NumberFormat
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.998'
$ws.Range("D39").NumberFormat = $__fmt
$ws.Range("E39").Value = '  +0.08%  '

$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$__fmt = $ws.Range("D40").NumberFormat
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.134'
$ws.Range("D40").NumberFormat = $__fmt
$ws.Range("E40").Value = '  -0.86%  '

$ws.Range("B41").Value = 'TheGraph'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$__fmt = $ws.Range("D41").NumberFormat
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.321'
$ws.Range("D41").NumberFormat = $__fmt
$ws.Range("E41").Value = '  -0.32%  '

$ws.Range("B42").Value = 'LidoDAOToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$__fmt = $ws.Range("D42").NumberFormat
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.33'
$ws.Range("D42").NumberFormat = $__fmt
$ws.Range("E42").Value = '  -1.34%  '

$__fmt = $ws.Range("D43").NumberFormat
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.70'
$ws.Range("D43").NumberFormat = $__fmt

$__fmt = $ws.Range("D44").NumberFormat
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '146.44'
$ws.Range("D44").NumberFormat = $__fmt
$ws.Range("E44").Value = '  +2.63%  '

$__fmt = $ws.Range("D45").NumberFormat
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.07'
$ws.Range("D45").NumberFormat = $__fmt
$ws.Range("E45").Value = '  +3.74%  '

$__fmt = $ws.Range("D46").NumberFormat
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.33'
$ws.Range("D46").NumberFormat = $__fmt
$ws.Range("E46").Value = '  +1.39%  '

$ws.Range("E47").Value = '  +12.38%  '

$ws.Range("D48").Value = '0.0₃0561'
$ws.Range("E48").Value = '  +29.50%  '

$ws.Range("E49").Value = '  -2.26%  '

$__fmt = $ws.Range("D50").NumberFormat
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.12'
$ws.Range("D50").NumberFormat = $__fmt
$ws.Range("E50").Value = '  -0.91%  '

$ws.Range("E51").Value = '  +1.06%  '

